$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
}

Set-TextValue "D2" "247.07"

Set-TextValue "D4" "5.493"

Set-TextValue "D5" "0.05630"

Set-TextValue "D6" "3.386"

Set-TextValue "D7" "6.481"

Set-TextValue "D8" "0.8039"

Set-TextValue "D9" "1.047"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1423"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07273"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03189"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.02948"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09264"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001661"
$ws.Range("E15").Value = "14BitForexTokenBF"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.242"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04704"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005910"
$ws.Range("E18").Value = "17OneONE"

Set-TextValue "D19" "0.006270"

Set-TextValue "D20" "0.001052"

Set-TextValue "D21" "0.004029"

Set-TextValue "D22" "0.0001503"

Set-TextValue "D23" "0.0003306"

Set-TextValue "D24" "3.971"

Set-TextValue "D25" "2.131"

$ws.Range("E27").Value = "26ProBitTokenPROBBestin24h"

Set-TextValue "D40" "0.04177"

$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1043"
$ws.Range("E41").Value = "40BKEXTokenBKK"

Set-TextValue "D42" "0.002975"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003248"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

Set-TextValue "D44" "0.009270"

Set-TextValue "D47" "0.6811"

Set-TextValue "D48" "0.02507"
$ws.Range("E48").Value = "47BOLOBOLO"

Set-TextValue "D50" "0.01012"
